$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each of the 4 data blocks (rows 2-9, 10-17, 18-25, 26-33) lists items 1-8.
# A new "item 9" summary row (49, 51, 53, 54, 52, 50) is inserted right
# after every block. Insert bottom-to-top so earlier row numbers stay valid.
$insertAt = 34, 26, 18, 10

foreach ($r in $insertAt) {
    $ws.Rows.Item($r).Insert()

    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = 49
    $ws.Cells.Item($r, 3).Value = 51
    $ws.Cells.Item($r, 4).Value = 53
    $ws.Cells.Item($r, 5).Value = 54
    $ws.Cells.Item($r, 6).Value = 52
    $ws.Cells.Item($r, 7).Value = 50

    # Drop the border formatting the insert inherited so the new row matches
    # the plain (unbordered) data rows, then restore column A's normal style.
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 7)).ClearFormats()
    $ws.Cells.Item(2, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
